$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (Excel would otherwise
# auto-convert plain numeric-looking strings like "619.36" into numbers).
# The cell's original style is restored afterwards so no stray number format
# is left applied to the cell.
function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "69.874.36"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").Value = "3.760.92"
$ws.Range("E3").Value = "  +3.33%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "619.36"
$ws.Range("E5").Value = "  +3.99%  "

Set-TextValue $ws.Range("D6") "177.63"
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("D7").Value = "3.760.88"
$ws.Range("E7").Value = "  +3.42%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -0.27%  "

Set-TextValue $ws.Range("D10") "0.168"
$ws.Range("E10").Value = "  +3.51%  "

$ws.Range("E11").Value = "  -5.29%  "

Set-TextValue $ws.Range("D12") "0.488"
$ws.Range("E12").Value = "  -1.94%  "

Set-TextValue $ws.Range("D13") "40.90"
$ws.Range("E13").Value = "  +0.72%  "

Set-TextValue $ws.Range("D14") "0.0000258"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "4.378.95"
$ws.Range("E15").Value = "  +3.09%  "

$ws.Range("D16").Value = "3.757.17"
$ws.Range("E16").Value = "  +3.00%  "

$ws.Range("D17").Value = "69.906.57"
$ws.Range("E17").Value = "  -1.44%  "

$ws.Range("E18").Value = "  +0.19%  "

Set-TextValue $ws.Range("D19") "7.54"
$ws.Range("E19").Value = "  +0.87%  "

Set-TextValue $ws.Range("D20") "509.19"
$ws.Range("E20").Value = "  -0.75%  "

Set-TextValue $ws.Range("D21") "16.59"
$ws.Range("E21").Value = "  -2.14%  "

Set-TextValue $ws.Range("D22") "9.62"
$ws.Range("E22").Value = "  +5.31%  "

Set-TextValue $ws.Range("D23") "0.723"
$ws.Range("E23").Value = "  -2.35%  "

Set-TextValue $ws.Range("D24") "2.51"
$ws.Range("E24").Value = "  +1.00%  "

Set-TextValue $ws.Range("D25") "86.79"
$ws.Range("E25").Value = "  -0.57%  "

Set-TextValue $ws.Range("D26") "13.12"
$ws.Range("E26").Value = "  -2.63%  "

Set-TextValue $ws.Range("D27") "11.03"
$ws.Range("E27").Value = "  +0.44%  "

Set-TextValue $ws.Range("D28") "0.0000135"
$ws.Range("E28").Value = "  +22.38%  "

Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.18%  "

Set-TextValue $ws.Range("D30") "2.50"
$ws.Range("E30").Value = "  -1.10%  "

$ws.Range("E31").Value = "  +4.68%  "

Set-TextValue $ws.Range("D32") "7.86"
$ws.Range("E32").Value = "  -3.74%  "

Set-TextValue $ws.Range("D33") "31.03"
$ws.Range("E33").Value = "  -1.25%  "

$ws.Range("E34").Value = "  -1.55%  "

$ws.Range("E35").Value = "  +0.01%  "

Set-TextValue $ws.Range("D36") "1.07"
$ws.Range("E36").Value = "  +5.93%  "

Set-TextValue $ws.Range("D37") "6.15"
$ws.Range("E37").Value = "  +0.77%  "

Set-TextValue $ws.Range("D38") "0.335"
$ws.Range("E38").Value = "  -3.08%  "

Set-TextValue $ws.Range("D39") "0.132"
$ws.Range("E39").Value = "  +2.27%  "

Set-TextValue $ws.Range("D40") "2.12"
$ws.Range("E40").Value = "  -1.80%  "

Set-TextValue $ws.Range("D41") "50.57"
$ws.Range("E41").Value = "  -0.67%  "

Set-TextValue $ws.Range("D42") "45.21"
$ws.Range("E42").Value = "  +0.06%  "

Set-TextValue $ws.Range("D43") "422.94"
$ws.Range("E43").Value = "  +1.90%  "

Set-TextValue $ws.Range("D44") "8.67"
$ws.Range("E44").Value = "  -1.49%  "

$ws.Range("D45").Value = "3.013.22"
$ws.Range("E45").Value = "  -3.50%  "

Set-TextValue $ws.Range("D46") "2.77"
$ws.Range("E46").Value = "  -1.38%  "

Set-TextValue $ws.Range("D47") "0.0361"
$ws.Range("E47").Value = "  -1.75%  "

Set-TextValue $ws.Range("D48") "27.30"
$ws.Range("E48").Value = "  -3.95%  "

Set-TextValue $ws.Range("D50") "138.43"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("E51").Value = "  -0.17%  "
